$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Formula = "=A$1+1"
$ws.Range("C1").Value = "test"

$ws.Rows.Item(1).RowHeight = 16

$fnt = $ws.Range("B1").Font
$fnt.Bold = $true
$fnt.Color = 255
$fnt.Name = "Aptos Narrow"

$ws.Range("B1").Select() | Out-Null
